$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 1 - header row (A1:O1), all using the same header style as A1 (s=3)
# ---------------------------------------------------------------------------
$headers = @(
    "Order Received Data and Time",
    "OrderID",
    "Emp ID-Order Assigned",
    "Assignee_QA",
    "Typist",
    "Typist QC",
    "Client",
    "Lob",
    "Process",
    "Product Name",
    "State",
    "County",
    "Municipality",
    "Status",
    "Tier"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
# Extend the header style (font/fill/border) from A1 onto the newly used E1:O1 cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row 2
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 45509.0625
$ws.Cells.Item(2, 2).Value = "FS18-001"
$ws.Cells.Item(2, 3).Value = "SIPL5316"
$ws.Cells.Item(2, 4).Value = "SIPL5688"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = "FAMS"
$ws.Cells.Item(2, 8).Value = "Servicing"
$ws.Cells.Item(2, 9).Value = "Search"
$ws.Cells.Item(2, 10).Value = "1 Owner - FCL Info"
$ws.Cells.Item(2, 11).Value = "AL"
$ws.Cells.Item(2, 12).Value = "Autauga"
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(2, 14).Value = "WIP"
$ws.Cells.Item(2, 15).Value = "Search(T1)"

# ---------------------------------------------------------------------------
# 3. Row 3
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = 45509.0625
$ws.Cells.Item(3, 2).Value = "FS18-001"
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "SIPL0102"
$ws.Cells.Item(3, 6).Value = "SIPL0103"
$ws.Cells.Item(3, 7).Value = "FAMS"
$ws.Cells.Item(3, 8).Value = "Servicing"
$ws.Cells.Item(3, 9).Value = "Typing"
$ws.Cells.Item(3, 10).Value = "1 Owner - FCL Info Update"
$ws.Cells.Item(3, 11).Value = "AL"
$ws.Cells.Item(3, 12).Value = "Baldwin"
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(3, 14).Value = "Typing"
$ws.Cells.Item(3, 15).Value = "Typing(T1)"

# Extend the existing bordered "data" style (copy from an already-bordered cell,
# e.g. C2 which uses style s=1) onto the brand-new M column (13) cells.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("M2:M3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Special styles
# ---------------------------------------------------------------------------
# K2/L2/K3/L3 -> centered, 10pt font, bordered (new cellXfs entry, fontId 19)
$c = $ws.Cells.Item(2, 11)
$c.Font.Size = 10
$c.HorizontalAlignment = -4108   # xlCenter
$c.VerticalAlignment = -4108     # xlCenter
$ws.Range("K2").Copy() | Out-Null
$ws.Range("L2:L2,K3:L3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# E3 -> font with explicit black color (new cellXfs entry, fontId 20, border 10)
$e3 = $ws.Cells.Item(3, 5)
$e3.Font.Color = 0

# F3 -> same font as E3, but without the left border (new cellXfs entry, fontId 20, border 11)
$f3 = $ws.Cells.Item(3, 6)
$f3.Font.Color = 0
$f3.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> xlLineStyleNone

# ---------------------------------------------------------------------------
# 5. Column widths (best effort match of the target stored widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666    # -> 12
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666    # -> 12
$ws.Columns.Item(5).ColumnWidth = 19.5                  # -> ~20.33203125
$ws.Columns.Item(6).ColumnWidth = 11.833333333333334    # -> ~12.6640625
$ws.Columns.Item(9).ColumnWidth = 11.833333333333334    # -> ~12.6640625
$ws.Columns.Item(10).ColumnWidth = 15.166666666666666   # -> 16
$ws.Columns.Item(11).ColumnWidth = 15.166666666666666   # -> 16
$ws.Columns.Item(12).ColumnWidth = 31.166666666666668   # -> 32
$ws.Columns.Item(15).ColumnWidth = 10.0                 # -> ~10.88671875
$ws.Columns.Item(16).ColumnWidth = 10.666666666666666   # -> ~11.5546875

# ---------------------------------------------------------------------------
# 6. Selection (matches the saved sheet view in the target file)
# ---------------------------------------------------------------------------
$ws.Range("D6").Select() | Out-Null

Write-Host "Edit complete"
